$wb = $excel.ActiveWorkbook

# ALC row 18 (Leve Item ID 5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4022.4856
$ws.Range("I18").Value = 4022.4856
$ws.Range("K18").Value = 4022.4856
$ws.Range("M18").Value = -3738.4856

# ALC row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 209
$ws.Range("I33").Value = 206.25
$ws.Range("J33").Value = 216.33333
$ws.Range("K33").Value = 206.25
$ws.Range("L33").Value = 216.33333
$ws.Range("M33").Value = 22.75
$ws.Range("N33").Value = -674.3333299999999

# ALC row 43 (Leve Item ID 5472)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5919.7144
$ws.Range("I43").Value = 5999.5
$ws.Range("K43").Value = 5999.5
$ws.Range("M43").Value = -5930.5

# ALC row 51 (Leve Item ID 5486)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5476.4
$ws.Range("I51").Value = 4349.25
$ws.Range("K51").Value = 4349.25
$ws.Range("M51").Value = -3865.25

# ALC row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4300
$ws.Range("I76").Value = 3600
$ws.Range("K76").Value = 3600
$ws.Range("M76").Value = -3285

# ALC row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4300
$ws.Range("I79").Value = 3600
$ws.Range("K79").Value = 3600
$ws.Range("M79").Value = -2508

# ALC row 112 (Leve Item ID 27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1571.5
$ws.Range("J112").Value = 2369.25
$ws.Range("L112").Value = 7107.75
$ws.Range("N112").Value = -9323.75

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4036.862
$ws.Range("I132").Value = 4058.6924
$ws.Range("K132").Value = 12176.0772
$ws.Range("M132").Value = -9646.0772

# ALC row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7905.3145
$ws.Range("I137").Value = 14531.4375
$ws.Range("K137").Value = 43594.3125
$ws.Range("M137").Value = -41044.3125

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3750.6885
$ws.Range("J138").Value = 4495.8774
$ws.Range("L138").Value = 13487.6322
$ws.Range("N138").Value = -23767.6322

# ARM row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5226.3447
$ws.Range("I102").Value = 5167.56
$ws.Range("J102").Value = 5593.75
$ws.Range("K102").Value = 5167.56
$ws.Range("L102").Value = 5593.75
$ws.Range("M102").Value = -3545.56
$ws.Range("N102").Value = -8837.75

# ARM row 125 (Leve Item ID 34251)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 159165.83
$ws.Range("J125").Value = 159165.83
$ws.Range("L125").Value = 159165.83
$ws.Range("N125").Value = -169005.83

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3920.75
$ws.Range("I132").Value = 2844.2173
$ws.Range("K132").Value = 8532.651899999999
$ws.Range("M132").Value = -6002.651899999999

# BSM row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3626.7
$ws.Range("I105").Value = 2323.8572
$ws.Range("K105").Value = 2323.8572
$ws.Range("M105").Value = -576.8571999999999

# CRP row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1044.9412
$ws.Range("I16").Value = 947.4666999999999
$ws.Range("K16").Value = 947.4666999999999
$ws.Range("M16").Value = -660.4666999999999

# CRP row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 999.8333
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 999.8
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 999.8
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -1699.8

# CRP row 86 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11937
$ws.Range("J86").Value = 13199.7
$ws.Range("L86").Value = 13199.7
$ws.Range("N86").Value = -15445.7

# CRP row 89 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11937
$ws.Range("J89").Value = 13199.7
$ws.Range("L89").Value = 65998.5
$ws.Range("N89").Value = -77230.5

# CRP row 107 (Leve Item ID 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 37047760
$ws.Range("I107").Value = 50013370
$ws.Range("K107").Value = 50013370
$ws.Range("M107").Value = -50011450

# CRP row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1044.9412
$ws.Range("I113").Value = 947.4666999999999
$ws.Range("K113").Value = 947.4666999999999
$ws.Range("M113").Value = 1222.5333

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5035.269
$ws.Range("I134").Value = 4165.478
$ws.Range("K134").Value = 12496.434
$ws.Range("M134").Value = -9961.434000000001

# CUL row 7 (Leve Item ID 4728)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 10
$ws.Range("K7").Value = 30
$ws.Range("M7").Value = 82

# CUL row 34 (Leve Item ID 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1564540.5
$ws.Range("J34").Value = 2842.5715
$ws.Range("L34").Value = 8527.7145
$ws.Range("N34").Value = -8695.7145

# CUL row 104 (Leve Item ID 19807)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 10000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 30000
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -35242

# CUL row 117 (Leve Item ID 27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3202.1667
$ws.Range("J117").Value = 1740.3334
$ws.Range("L117").Value = 5221.0002
$ws.Range("N117").Value = -12105.0002

# CUL row 132 (Leve Item ID 43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 101995.4
$ws.Range("I132").Value = 988.5
$ws.Range("K132").Value = 8896.5
$ws.Range("M132").Value = -6366.5

# GSM row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9685.143
$ws.Range("I70").Value = 6699.5
$ws.Range("K70").Value = 6699.5
$ws.Range("M70").Value = -6429.5

# GSM row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9685.143
$ws.Range("I73").Value = 6699.5
$ws.Range("K73").Value = 6699.5
$ws.Range("M73").Value = -5763.5

# GSM row 99 (Leve Item ID 19532)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 19234.75
$ws.Range("I99").Value = 15646.667
$ws.Range("K99").Value = 15646.667
$ws.Range("M99").Value = -13400.667

# GSM row 123 (Leve Item ID 34150)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 21899.8
$ws.Range("J123").Value = 21899.8
$ws.Range("L123").Value = 21899.8
$ws.Range("N123").Value = -26799.8

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3708.3333
$ws.Range("I132").Value = 3967.5454
$ws.Range("K132").Value = 11902.6362
$ws.Range("M132").Value = -9372.636200000001

# LTW row 22 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2229.9167
$ws.Range("I22").Value = 2487.8667
$ws.Range("J22").Value = 1800
$ws.Range("K22").Value = 2487.8667
$ws.Range("L22").Value = 1800
$ws.Range("M22").Value = -2192.8667
$ws.Range("N22").Value = -2390

# LTW row 27 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2229.9167
$ws.Range("I27").Value = 2487.8667
$ws.Range("J27").Value = 1800
$ws.Range("K27").Value = 2487.8667
$ws.Range("L27").Value = 1800
$ws.Range("M27").Value = -2380.8667
$ws.Range("N27").Value = -2014

# LTW row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 702.94116
$ws.Range("I55").Value = 708.3333
$ws.Range("K55").Value = 708.3333
$ws.Range("M55").Value = -535.3333

# LTW row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4001.5
$ws.Range("I122").Value = 3632.7
$ws.Range("K122").Value = 10898.1
$ws.Range("M122").Value = -8448.099999999999

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10035.8125
$ws.Range("J136").Value = 11054.777
$ws.Range("L136").Value = 33164.331
$ws.Range("N136").Value = -38264.331

# WVR row 74 (Leve Item ID 19022)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 55359.5
$ws.Range("I74").Value = 11545.667
$ws.Range("K74").Value = 11545.667
$ws.Range("M74").Value = -10609.667

# WVR row 77 (Leve Item ID 19022)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 55359.5
$ws.Range("I77").Value = 11545.667
$ws.Range("K77").Value = 34637.001
$ws.Range("M77").Value = -29957.001

# WVR row 96 (Leve Item ID 19977)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1965.5883
$ws.Range("I96").Value = 1280.8334
$ws.Range("K96").Value = 1280.8334
$ws.Range("M96").Value = 92.16660000000002

# WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 25217.37
$ws.Range("I126").Value = 35778.75
$ws.Range("K126").Value = 107336.25
$ws.Range("M126").Value = -104866.25

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3287.2307
$ws.Range("I136").Value = 2441.8572
$ws.Range("K136").Value = 7325.571599999999
$ws.Range("M136").Value = -4775.571599999999
